$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row - Right column (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row - Right column (B12): 27 -> 45
$ws.Range("B12").Value = 45

# Update "Total" row - Max column (E12): correct/total marks text "25/84" -> "45/140"
$ws.Range("E12").Value = "45/140"
